$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (labels), with two extra "bold" summary columns inserted
# right after the leading blank cell, and the country columns re-ordered /
# renamed to match the new dataset layout.
# (A1 is left untouched - it already holds the empty-string label.)
$ws.Range("B1").Value = "`$ bold(All)"
$ws.Range("C1").Value = "`$ bold(Europe)"
$ws.Range("D1").Value = "France"
$ws.Range("E1").Value = "Germany"
$ws.Range("F1").Value = "Italy"
$ws.Range("G1").Value = "Poland"
$ws.Range("H1").Value = "Spain"
$ws.Range("I1").Value = "United Kingdom"
$ws.Range("J1").Value = "Switzerland"
$ws.Range("K1").Value = "Japan"
$ws.Range("L1").Value = "Saudi Arabia"
$ws.Range("M1").Value = "USA"

# Row 2 values (A2's label text is unchanged by the edit, only its shared-
# string index shifts because of the two new strings inserted earlier in
# the table - re-assert the original text, CRLF included, verbatim)
$ws.Range("A2").Value = "Belief about GCS support in the U.S.`r`n(except for the U.S.: support in the EU)"
$ws.Range("B2").Value = 34.3455077975836
$ws.Range("C2").Value = 28.1230920639193
$ws.Range("D2").Value = 28.0327920637749
$ws.Range("E2").Value = 24.2540573483137
$ws.Range("F2").Value = 29.499120537839
$ws.Range("G2").Value = 32.294256431363
$ws.Range("H2").Value = 29.3982179553451
$ws.Range("I2").Value = 28.8591089308659
$ws.Range("J2").Value = 25.3832571311288
$ws.Range("K2").Value = 30.7964256295476
$ws.Range("L2").Value = 45.291373979237
$ws.Range("M2").Value = 42.4040946837611

# Row 3 values
$ws.Range("A3").Value = "Belief about GCS support in own country"
$ws.Range("B3").Value = 42.5574333274623
$ws.Range("C3").Value = 42.880295347068
$ws.Range("D3").Value = 44.6898148665286
$ws.Range("E3").Value = 38.8014169530983
$ws.Range("F3").Value = 44.370617045058
$ws.Range("G3").Value = 41.7048669024782
$ws.Range("H3").Value = 43.9428615425684
$ws.Range("I3").Value = 44.8409655646873
$ws.Range("J3").Value = 43.8248904709614
$ws.Range("K3").Value = 37.2281838448364
$ws.Range("L3").Value = 54.9042600841208
$ws.Range("M3").Value = 43.1563898238002
